# Performance Metrics - Garrett PC.xlsx
# "Ran Perf-Met on VectorAdd. Updated spreadsheet"
#
# Changes:
#   1. Header E1 "Page-faults" -> "Page-faults (avg)"
#   2. New results row (row 4, VectorAdd) filled in with B4:G4 raw metrics
#      and H4 = F4/G4 (CPI-style ratio formula)
#   3. Column E gets an explicit width (new metric column now has data)
#   4. Active cell / selection moves from B4 to H5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Page-faults" column header (column E) to "Page-faults (avg)"
$ws.Range("E1").Value2 = "Page-faults (avg)"

# 2. Fill in the VectorAdd performance numbers on row 4
$ws.Range("B4").Value2 = 0.002838102
$ws.Range("C4").Value2 = 274478.7
$ws.Range("D4").Value2 = 64397.7
$ws.Range("E4").Value2 = 104.7
$ws.Range("F4").Value2 = 27478802.9
$ws.Range("G4").Value2 = 33095732.6
$ws.Range("H4").Formula = "=F4/G4"

# 3. Give column E an explicit width (closest reachable value to 16.53
#    characters under Excel's pixel-quantized column-width model)
$ws.Columns("E").ColumnWidth = 15.6

# 4. Move the active selection to H5
$ws.Range("H5").Select() | Out-Null
